$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 566 (pushes the existing row 566.."657" down
# by one, growing the used range from A1:R657 to A1:R658).
$ws.Rows.Item(566).Insert()

# Populate the freshly inserted row with the new weekly price observation.
$ws.Range("A566").Value = 3
$ws.Range("B566").Value = "Femacal de La Calera"
$ws.Range("C566").Value = "Coquimbo"
$ws.Range("D566").Value = 45218
$ws.Range("E566").Value = 5
$ws.Range("F566").Value = 100112040
$ws.Range("G566").Value = "Cilantro"
$ws.Range("H566").Value = "Sin especificar"
$ws.Range("I566").Value = "Primera"
$ws.Range("J566").Value = 173
$ws.Range("K566").Value = 3500
$ws.Range("L566").Value = 3800
$ws.Range("M566").Value = 3630
$ws.Range("N566").Value = "$/docena de atados (3 kilos)"
$ws.Range("O566").Value = "Provincia de Quillota"
$ws.Range("P566").Value = 1210
$ws.Range("Q566").Value = 3
$ws.Range("R566").Value = "Hortaliza"
